$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Lgi1"
$ws.Range("C2").Value = "Adam23"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.000409
$ws.Range("H2").Value = 0.001227
$ws.Range("I2").Value = 0.164521319388576
$ws.Range("J2").Value = 0.228024530756365
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.09286699999999999
$ws.Range("N2").Value = 0.278601
$ws.Range("O2").Value = 0.003009076821730935
$ws.Range("P2").Value = 0.003071957783644885
$ws.Range("Q2").Value = 0.000037982603
$ws.Range("R2").Value = 0.000341843427
$ws.Range("S2").Value = 0.0004950572888527564
$ws.Range("T2").Value = 0.000700481732118988

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Lgi1"
$ws.Range("C3").Value = "Adam23"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.000409
$ws.Range("H3").Value = 0.001227
$ws.Range("I3").Value = 0.164521319388576
$ws.Range("J3").Value = 0.228024530756365
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 28.61718033333333
$ws.Range("N3").Value = 85.851541
$ws.Range("O3").Value = 0.9272539658256183
$ws.Range("P3").Value = 0.9466308793322996
$ws.Range("Q3").Value = 0.01170442675633333
$ws.Range("R3").Value = 0.105339840807
$ws.Range("S3").Value = 0.1525530458659203
$ws.Range("T3").Value = 0.2158550620592328

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Lgi1"
$ws.Range("C4").Value = "Adam23"
$ws.Range("D4").Value = "Inflammatory-Mac"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.000409
$ws.Range("H4").Value = 0.001227
$ws.Range("I4").Value = 0.164521319388576
$ws.Range("J4").Value = 0.228024530756365
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.243062
$ws.Range("N4").Value = 0.729186
$ws.Range("O4").Value = 0.007875695677081898
$ws.Range("P4").Value = 0.00804027483183793
$ws.Range("Q4").Value = 0.00009941235800000001
$ws.Range("R4").Value = 0.0008947112220000001
$ws.Range("S4").Value = 0.001295719843896419
$ws.Range("T4").Value = 0.001833379895682056

# Row 5
$ws.Range("A5").Value = "ECs"
$ws.Range("B5").Value = "Lgi1"
$ws.Range("C5").Value = "Adam23"
$ws.Range("D5").Value = "MuSCs"
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.000409
$ws.Range("H5").Value = 0.001227
$ws.Range("I5").Value = 0.164521319388576
$ws.Range("J5").Value = 0.228024530756365
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 1.8951925
$ws.Range("N5").Value = 3.790385
$ws.Range("O5").Value = 0.06140803325689756
$ws.Range("P5").Value = 0.04179418847656979
$ws.Range("Q5").Value = 0.0007751337325000001
$ws.Range("R5").Value = 0.004650802395
$ws.Range("S5").Value = 0.01010293065248234
$ws.Range("T5").Value = 0.009530100215712904

# Row 6
$ws.Range("A6").Value = "ECs"
$ws.Range("B6").Value = "Lgi1"
$ws.Range("C6").Value = "Adam23"
$ws.Range("D6").Value = "Resolving-Mac"
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.3333333333333333
$ws.Range("G6").Value = 0.000409
$ws.Range("H6").Value = 0.001227
$ws.Range("I6").Value = 0.164521319388576
$ws.Range("J6").Value = 0.228024530756365
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.3333333333333333
$ws.Range("M6").Value = 0.01398766666666667
$ws.Range("N6").Value = 0.041963
$ws.Range("O6").Value = 0.000453228418671488
$ws.Range("P6").Value = 0.0004626995756479349
$ws.Range("Q6").Value = 0.000005720955666666667
$ws.Range("R6").Value = 0.000051488601
$ws.Range("S6").Value = 0.00007456573742423114
$ws.Range("T6").Value = 0.0001055068536182896

# Row 7
$ws.Range("A7").Value = "MuSCs"
$ws.Range("B7").Value = "Lgi1"
$ws.Range("C7").Value = "Adam23"
$ws.Range("D7").Value = "ECs"
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 0.5
$ws.Range("G7").Value = 0.002077
$ws.Range("H7").Value = 0.004154
$ws.Range("I7").Value = 0.8354786806114239
$ws.Range("J7").Value = 0.771975469243635
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 0.6666666666666666
$ws.Range("M7").Value = 0.09286699999999999
$ws.Range("N7").Value = 0.278601
$ws.Range("O7").Value = 0.003009076821730935
$ws.Range("P7").Value = 0.003071957783644885
$ws.Range("Q7").Value = 0.000192884759
$ws.Range("R7").Value = 0.001157308554
$ws.Range("S7").Value = 0.002514019532878178
$ws.Range("T7").Value = 0.002371476051525897

# Row 8
$ws.Range("A8").Value = "MuSCs"
$ws.Range("B8").Value = "Lgi1"
$ws.Range("C8").Value = "Adam23"
$ws.Range("D8").Value = "FAPs"
$ws.Range("E8").Value = 1
$ws.Range("F8").Value = 0.5
$ws.Range("G8").Value = 0.002077
$ws.Range("H8").Value = 0.004154
$ws.Range("I8").Value = 0.8354786806114239
$ws.Range("J8").Value = 0.771975469243635
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 28.61718033333333
$ws.Range("N8").Value = 85.851541
$ws.Range("O8").Value = 0.9272539658256183
$ws.Range("P8").Value = 0.9466308793322996
$ws.Range("Q8").Value = 0.05943788355233333
$ws.Range("R8").Value = 0.3566273013139999
$ws.Range("S8").Value = 0.7747009199596979
$ws.Range("T8").Value = 0.7307758172730668

# Row 9
$ws.Range("A9").Value = "MuSCs"
$ws.Range("B9").Value = "Lgi1"
$ws.Range("C9").Value = "Adam23"
$ws.Range("D9").Value = "Inflammatory-Mac"
$ws.Range("E9").Value = 1
$ws.Range("F9").Value = 0.5
$ws.Range("G9").Value = 0.002077
$ws.Range("H9").Value = 0.004154
$ws.Range("I9").Value = 0.8354786806114239
$ws.Range("J9").Value = 0.771975469243635
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.243062
$ws.Range("N9").Value = 0.729186
$ws.Range("O9").Value = 0.007875695677081898
$ws.Range("P9").Value = 0.00804027483183793
$ws.Range("Q9").Value = 0.0005048397739999999
$ws.Range("R9").Value = 0.003029038644
$ws.Range("S9").Value = 0.006579975833185479
$ws.Range("T9").Value = 0.006206894936155875

# Row 10
$ws.Range("A10").Value = "MuSCs"
$ws.Range("B10").Value = "Lgi1"
$ws.Range("C10").Value = "Adam23"
$ws.Range("D10").Value = "MuSCs"
$ws.Range("E10").Value = 1
$ws.Range("F10").Value = 0.5
$ws.Range("G10").Value = 0.002077
$ws.Range("H10").Value = 0.004154
$ws.Range("I10").Value = 0.8354786806114239
$ws.Range("J10").Value = 0.771975469243635
$ws.Range("K10").Value = 2
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 1.8951925
$ws.Range("N10").Value = 3.790385
$ws.Range("O10").Value = 0.06140803325689756
$ws.Range("P10").Value = 0.04179418847656979
$ws.Range("Q10").Value = 0.0039363148225
$ws.Range("R10").Value = 0.01574525929
$ws.Range("S10").Value = 0.05130510260441522
$ws.Range("T10").Value = 0.03226408826085689

# Row 11
$ws.Range("A11").Value = "MuSCs"
$ws.Range("B11").Value = "Lgi1"
$ws.Range("C11").Value = "Adam23"
$ws.Range("D11").Value = "Resolving-Mac"
$ws.Range("E11").Value = 1
$ws.Range("F11").Value = 0.5
$ws.Range("G11").Value = 0.002077
$ws.Range("H11").Value = 0.004154
$ws.Range("I11").Value = 0.8354786806114239
$ws.Range("J11").Value = 0.771975469243635
$ws.Range("K11").Value = 1
$ws.Range("L11").Value = 0.3333333333333333
$ws.Range("M11").Value = 0.01398766666666667
$ws.Range("N11").Value = 0.041963
$ws.Range("O11").Value = 0.000453228418671488
$ws.Range("P11").Value = 0.0004626995756479349
$ws.Range("Q11").Value = 0.00002905238366666667
$ws.Range("R11").Value = 0.000174314302
$ws.Range("S11").Value = 0.0003786626812472569
$ws.Range("T11").Value = 0.0003571927220296454

